# ContactUs sheet update: the single seeded contact's phone/email placeholder
# data is replaced with a real verified contact, and the ad-hoc per-row
# "mailto:" hyperlinks are removed now that every row shares one address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all "mailto:" hyperlinks from the Email column. The cells keep
# whatever formatting they already had (the blue/underlined hyperlink look).
$ws.Hyperlinks.Delete()

# Re-key PhoneNumber (column F) as quote-prefixed text (was a plain number)
# before touching Email so the shared-string table picks up "9876543210"
# ahead of the new email address, matching the authored edit order.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Value = "'9876543210"
}

# Every contact row now points at the same updated email address.
$ws.Range("E2:E6").Value = "Prashant.Chandra@ascendlearning.com"

# Widen the Email column so the longer address fits without truncation.
$ws.Columns("E").ColumnWidth = 36.1666666

# Leave the selection where the author ended up after editing.
$ws.Range("E8").Select()
